# Applies the cryptos-list price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = Price, Column E = Volume(1h). Price values that look like plain
# numbers (e.g. "29.10", "1.00") are entered with a leading apostrophe so Excel
# keeps them as literal text (matching the target formatting) instead of
# collapsing them to numeric values and losing trailing zeros.

# Row 2: D2='67.103.05', E2='  +0.70%  '
$ws.Cells.Item(2, 4).Value = '67.103.05'
$ws.Cells.Item(2, 5).Value = '  +0.70%  '

# Row 3: D3='3.509.57', E3='  +0.15%  '
$ws.Cells.Item(3, 4).Value = '3.509.57'
$ws.Cells.Item(3, 5).Value = '  +0.15%  '

# Row 4: E4='  +0.00%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

# Row 5: D5='594.96', E5='  +0.58%  '
$ws.Cells.Item(5, 4).Value = '''594.96'
$ws.Cells.Item(5, 5).Value = '  +0.58%  '

# Row 6: D6='173.45', E6='  +2.10%  '
$ws.Cells.Item(6, 4).Value = '''173.45'
$ws.Cells.Item(6, 5).Value = '  +2.10%  '

# Row 7: D7='0.999', E7='  -0.02%  '
$ws.Cells.Item(7, 4).Value = '''0.999'
$ws.Cells.Item(7, 5).Value = '  -0.02%  '

# Row 8: E8='  +0.98%  '
$ws.Cells.Item(8, 5).Value = '  +0.98%  '

# Row 9: E9='  +4.80%  '
$ws.Cells.Item(9, 5).Value = '  +4.80%  '

# Row 10: E10='  -0.63%  '
$ws.Cells.Item(10, 5).Value = '  -0.63%  '

# Row 11: D11='0.435', E11='  -1.09%  '
$ws.Cells.Item(11, 4).Value = '''0.435'
$ws.Cells.Item(11, 5).Value = '  -1.09%  '

# Row 12: D12='4.114.95', E12='  +0.11%  '
$ws.Cells.Item(12, 4).Value = '4.114.95'
$ws.Cells.Item(12, 5).Value = '  +0.11%  '

# Row 13: E13='  +0.02%  '
$ws.Cells.Item(13, 5).Value = '  +0.02%  '

# Row 14: D14='29.10', E14='  +2.64%  '
$ws.Cells.Item(14, 4).Value = '''29.10'
$ws.Cells.Item(14, 5).Value = '  +2.64%  '

# Row 15: D15='67.067.68', E15='  +0.60%  '
$ws.Cells.Item(15, 4).Value = '67.067.68'
$ws.Cells.Item(15, 5).Value = '  +0.60%  '

# Row 16: E16='  +0.64%  '
$ws.Cells.Item(16, 5).Value = '  +0.64%  '

# Row 17: D17='3.499.95', E17='  -0.16%  '
$ws.Cells.Item(17, 4).Value = '3.499.95'
$ws.Cells.Item(17, 5).Value = '  -0.16%  '

# Row 19: D19='14.18', E19='  +0.59%  '
$ws.Cells.Item(19, 4).Value = '''14.18'
$ws.Cells.Item(19, 5).Value = '  +0.59%  '

# Row 20: D20='396.22', E20='  +1.25%  '
$ws.Cells.Item(20, 4).Value = '''396.22'
$ws.Cells.Item(20, 5).Value = '  +1.25%  '

# Row 21: E21='  +0.75%  '
$ws.Cells.Item(21, 5).Value = '  +0.75%  '

# Row 22: D22='73.16', E22='  +0.02%  '
$ws.Cells.Item(22, 4).Value = '''73.16'
$ws.Cells.Item(22, 5).Value = '  +0.02%  '

# Row 23: D23='1.00', E23='  -0.08%  '
$ws.Cells.Item(23, 4).Value = '''1.00'
$ws.Cells.Item(23, 5).Value = '  -0.08%  '

# Row 24: D24='0.538'
$ws.Cells.Item(24, 4).Value = '''0.538'

# Row 25: D25='5.70', E25='  -3.03%  '
$ws.Cells.Item(25, 4).Value = '''5.70'
$ws.Cells.Item(25, 5).Value = '  -3.03%  '

# Row 26: E26='  -0.87%  '
$ws.Cells.Item(26, 5).Value = '  -0.87%  '

# Row 27: D27='10.30', E27='  +0.41%  '
$ws.Cells.Item(27, 4).Value = '''10.30'
$ws.Cells.Item(27, 5).Value = '  +0.41%  '

# Row 28: E28='  +0.59%  '
$ws.Cells.Item(28, 5).Value = '  +0.59%  '

# Row 29: E29='  -0.15%  '
$ws.Cells.Item(29, 5).Value = '  -0.15%  '

# Row 30: D30='6.32', E30='  -0.78%  '
$ws.Cells.Item(30, 4).Value = '''6.32'
$ws.Cells.Item(30, 5).Value = '  -0.78%  '

# Row 31: E31='  -3.07%  '
$ws.Cells.Item(31, 5).Value = '  -3.07%  '

# Row 32: E32='  -0.17%  '
$ws.Cells.Item(32, 5).Value = '  -0.17%  '

# Row 33: D33='23.87', E33='  +1.08%  '
$ws.Cells.Item(33, 4).Value = '''23.87'
$ws.Cells.Item(33, 5).Value = '  +1.08%  '

# Row 34: E34='  -0.85%  '
$ws.Cells.Item(34, 5).Value = '  -0.85%  '

# Row 35: D35='1.68', E35='  +3.61%  '
$ws.Cells.Item(35, 4).Value = '''1.68'
$ws.Cells.Item(35, 5).Value = '  +3.61%  '

# Row 36: D36='163.29', E36='  +0.22%  '
$ws.Cells.Item(36, 4).Value = '''163.29'
$ws.Cells.Item(36, 5).Value = '  +0.22%  '

# Row 37: E37='  -0.04%  '
$ws.Cells.Item(37, 5).Value = '  -0.04%  '

# Row 38: E38='  -0.18%  '
$ws.Cells.Item(38, 5).Value = '  -0.18%  '

# Row 39: D39='7.08', E39='  +4.12%  '
$ws.Cells.Item(39, 4).Value = '''7.08'
$ws.Cells.Item(39, 5).Value = '  +4.12%  '

# Row 40: D40='4.70'
$ws.Cells.Item(40, 4).Value = '''4.70'

# Row 41: D41='0.0749', E41='  -0.09%  '
$ws.Cells.Item(41, 4).Value = '''0.0749'
$ws.Cells.Item(41, 5).Value = '  -0.09%  '

# Row 42: D42='27.53', E42='  +1.93%  '
$ws.Cells.Item(42, 4).Value = '''27.53'
$ws.Cells.Item(42, 5).Value = '  +1.93%  '

# Row 43: D43='26.47', E43='  -0.50%  '
$ws.Cells.Item(43, 4).Value = '''26.47'
$ws.Cells.Item(43, 5).Value = '  -0.50%  '

# Row 44: D44='2.815.10', E44='  -0.11%  '
$ws.Cells.Item(44, 4).Value = '2.815.10'
$ws.Cells.Item(44, 5).Value = '  -0.11%  '

# Row 45: D45='2.59', E45='  +2.69%  '
$ws.Cells.Item(45, 4).Value = '''2.59'
$ws.Cells.Item(45, 5).Value = '  +2.69%  '

# Row 46: E46='  -0.79%  '
$ws.Cells.Item(46, 5).Value = '  -0.79%  '

# Row 47: E47='  -2.51%  '
$ws.Cells.Item(47, 5).Value = '  -2.51%  '

# Row 48: D48='340.48', E48='  -4.07%  '
$ws.Cells.Item(48, 4).Value = '''340.48'
$ws.Cells.Item(48, 5).Value = '  -4.07%  '

# Row 49: D49='34.82', E49='  +3.35%  '
$ws.Cells.Item(49, 4).Value = '''34.82'
$ws.Cells.Item(49, 5).Value = '  +3.35%  '

# Row 50: E50='  -0.36%  '
$ws.Cells.Item(50, 5).Value = '  -0.36%  '

# Row 51: E51='  -0.86%  '
$ws.Cells.Item(51, 5).Value = '  -0.86%  '
